# Weekly data refresh: a new "Femacal de La Calera" / Perejil price
# observation (date serial 44574) was inserted above the existing block
# of rows, pushing the previously-recorded rows 19-23 down to 20-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 19; rows 19:23 (and their formatting) shift
# down to 20:24, matching the diff's row-shift pattern exactly.
$ws.Rows("19:19").Insert()

# Populate the new row 19 with the latest observation.
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = "Femacal de La Calera"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 44574
$ws.Range("E19").Value = 5
$ws.Range("F19").Value = 100112044
$ws.Range("G19").Value = "Perejil"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 3000
$ws.Range("N19").Value = "$/docena de atados (3 kilos)"
$ws.Range("O19").Value = "Provincia de Quillota"
$ws.Range("P19").Value = 1000
$ws.Range("Q19").Value = 3
$ws.Range("R19").Value = "Hortaliza"
